$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 15436.7
$ws.Range("I43").Value = 26580
$ws.Range("K43").Value = 26580
$ws.Range("M43").Value = -26511
$ws.Range("H70").Value = 103767.4
$ws.Range("J70").Value = 103767.4
$ws.Range("L70").Value = 311302.2
$ws.Range("N70").Value = -311842.2
$ws.Range("H73").Value = 103767.4
$ws.Range("J73").Value = 103767.4
$ws.Range("L73").Value = 311302.2
$ws.Range("N73").Value = -313174.2
$ws.Range("H111").Value = 16314.0625
$ws.Range("I111").Value = 15126.286
$ws.Range("K111").Value = 45378.858
$ws.Range("M111").Value = -42311.858
$ws.Range("H132").Value = 1456.5555
$ws.Range("I132").Value = 1212.4572
$ws.Range("K132").Value = 3637.3716
$ws.Range("M132").Value = -1107.3716
$ws.Range("H135").Value = 1018.8182
$ws.Range("J135").Value = 1399.5
$ws.Range("L135").Value = 12595.5
$ws.Range("N135").Value = -17665.5
$ws.Range("H137").Value = 1699.6
$ws.Range("I137").Value = 1500.4
$ws.Range("J137").Value = 1898.8
$ws.Range("K137").Value = 4501.200000000001
$ws.Range("L137").Value = 5696.4
$ws.Range("M137").Value = -1951.200000000001
$ws.Range("N137").Value = -10796.4
$ws.Range("H138").Value = 5816336.5
$ws.Range("I138").Value = 1006.1667
$ws.Range("K138").Value = 3018.5001
$ws.Range("M138").Value = 2121.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 5493.5713
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5493.5713
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5493.5713
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6131.5713
$ws.Range("H61").Value = 3375.014
$ws.Range("I61").Value = 2239.625
$ws.Range("J61").Value = 5645.7915
$ws.Range("K61").Value = 2239.625
$ws.Range("L61").Value = 5645.7915
$ws.Range("M61").Value = -2027.625
$ws.Range("N61").Value = -6069.7915
$ws.Range("H97").Value = 1485.3077
$ws.Range("I97").Value = 1053.238
$ws.Range("J97").Value = 3300
$ws.Range("K97").Value = 1053.238
$ws.Range("L97").Value = 3300
$ws.Range("M97").Value = -557.2380000000001
$ws.Range("N97").Value = -4292
$ws.Range("H135").Value = 40625
$ws.Range("J135").Value = 40625
$ws.Range("L135").Value = 40625
$ws.Range("N135").Value = -50765
$ws.Range("H136").Value = 3375.014
$ws.Range("I136").Value = 2239.625
$ws.Range("J136").Value = 5645.7915
$ws.Range("K136").Value = 6718.875
$ws.Range("L136").Value = 16937.3745
$ws.Range("M136").Value = -4168.875
$ws.Range("N136").Value = -22037.3745

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3222.6978
$ws.Range("I20").Value = 2766.1667
$ws.Range("J20").Value = 3799.3684
$ws.Range("K20").Value = 2766.1667
$ws.Range("L20").Value = 3799.3684
$ws.Range("M20").Value = -2519.1667
$ws.Range("N20").Value = -4293.368399999999
$ws.Range("H105").Value = 1482.4166
$ws.Range("I105").Value = 1548.8889
$ws.Range("J105").Value = 1283
$ws.Range("K105").Value = 1548.8889
$ws.Range("L105").Value = 1283
$ws.Range("M105").Value = 198.1111000000001
$ws.Range("N105").Value = -4777
$ws.Range("H122").Value = 40000
$ws.Range("J122").Value = 40000
$ws.Range("L122").Value = 40000
$ws.Range("N122").Value = -49800
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 1444.4193
$ws.Range("I134").Value = 1433.3077
$ws.Range("J134").Value = 1950
$ws.Range("K134").Value = 4299.9231
$ws.Range("L134").Value = 5850
$ws.Range("M134").Value = -1764.9231
$ws.Range("N134").Value = -10920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 350.70587
$ws.Range("I7").Value = 373.55554
$ws.Range("J7").Value = 325
$ws.Range("K7").Value = 373.55554
$ws.Range("L7").Value = 325
$ws.Range("M7").Value = -260.55554
$ws.Range("N7").Value = -551
$ws.Range("H31").Value = 58512.89
$ws.Range("I31").Value = 93497.63
$ws.Range("J31").Value = 3536.8572
$ws.Range("K31").Value = 93497.63
$ws.Range("L31").Value = 3536.8572
$ws.Range("M31").Value = -93202.63
$ws.Range("N31").Value = -4126.8572
$ws.Range("H34").Value = 58512.89
$ws.Range("I34").Value = 93497.63
$ws.Range("J34").Value = 3536.8572
$ws.Range("K34").Value = 93497.63
$ws.Range("L34").Value = 3536.8572
$ws.Range("M34").Value = -93295.63
$ws.Range("N34").Value = -3940.8572
$ws.Range("H132").Value = 2827.889
$ws.Range("I132").Value = 2896
$ws.Range("J132").Value = 2436.25
$ws.Range("K132").Value = 8688
$ws.Range("L132").Value = 7308.75
$ws.Range("M132").Value = -6158
$ws.Range("N132").Value = -12368.75
$ws.Range("H134").Value = 39658.9
$ws.Range("I134").Value = 21680
$ws.Range("J134").Value = 57637.8
$ws.Range("K134").Value = 65040
$ws.Range("L134").Value = 172913.4
$ws.Range("M134").Value = -62505
$ws.Range("N134").Value = -177983.4
$ws.Range("H135").Value = 80709
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 3919.1667
$ws.Range("I28").Value = 3919.1667
$ws.Range("K28").Value = 11757.5001
$ws.Range("M28").Value = -11525.5001
$ws.Range("H45").Value = 33335500
$ws.Range("J45").Value = 3250
$ws.Range("L45").Value = 9750
$ws.Range("N45").Value = -10814
$ws.Range("H98").Value = 529.4286
$ws.Range("I98").Value = 390
$ws.Range("J98").Value = 606.8889
$ws.Range("K98").Value = 1170
$ws.Range("L98").Value = 1820.6667
$ws.Range("M98").Value = 328
$ws.Range("N98").Value = -4816.6667
$ws.Range("H113").Value = 831.36
$ws.Range("J113").Value = 826.41174
$ws.Range("L113").Value = 2479.23522
$ws.Range("N113").Value = -6819.23522
$ws.Range("H114").Value = 2816.25
$ws.Range("I114").Value = 1633
$ws.Range("J114").Value = 3999.5
$ws.Range("K114").Value = 4899
$ws.Range("L114").Value = 11998.5
$ws.Range("M114").Value = -1645
$ws.Range("N114").Value = -18506.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 6666.5
$ws.Range("I57").Value = 6666.5
$ws.Range("K57").Value = 6666.5
$ws.Range("M57").Value = -5846.5
$ws.Range("H70").Value = 9680.913
$ws.Range("I70").Value = 9051.134
$ws.Range("J70").Value = 10861.75
$ws.Range("K70").Value = 9051.134
$ws.Range("L70").Value = 10861.75
$ws.Range("M70").Value = -8781.134
$ws.Range("N70").Value = -11401.75
$ws.Range("H73").Value = 9680.913
$ws.Range("I73").Value = 9051.134
$ws.Range("J73").Value = 10861.75
$ws.Range("K73").Value = 9051.134
$ws.Range("L73").Value = 10861.75
$ws.Range("M73").Value = -8115.134
$ws.Range("N73").Value = -12733.75
$ws.Range("H102").Value = 28573752
$ws.Range("I102").Value = 2287.6333
$ws.Range("J102").Value = 200002540
$ws.Range("K102").Value = 2287.6333
$ws.Range("L102").Value = 200002540
$ws.Range("M102").Value = -665.6333
$ws.Range("N102").Value = -200005784
$ws.Range("H113").Value = 4237.6665
$ws.Range("I113").Value = 3517.5
$ws.Range("J113").Value = 9999
$ws.Range("K113").Value = 3517.5
$ws.Range("L113").Value = 9999
$ws.Range("M113").Value = -1347.5
$ws.Range("N113").Value = -14339
$ws.Range("H122").Value = 2508.1143
$ws.Range("I122").Value = 2138.5
$ws.Range("K122").Value = 6415.5
$ws.Range("M122").Value = -3965.5
$ws.Range("H132").Value = 2943.5652
$ws.Range("J132").Value = 4500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560
$ws.Range("H140").Value = 74423.375
$ws.Range("J140").Value = 74423.375
$ws.Range("L140").Value = 74423.375
$ws.Range("N140").Value = -84783.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 81988.30499999999
$ws.Range("I61").Value = 88237.336
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 88237.336
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -88035.336
$ws.Range("N61").Value = -7404
$ws.Range("H113").Value = 81988.30499999999
$ws.Range("I113").Value = 88237.336
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 88237.336
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -86067.336
$ws.Range("N113").Value = -11340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H70").Value = 39052.5
$ws.Range("I70").Value = 38000
$ws.Range("K70").Value = 38000
$ws.Range("M70").Value = -37685
$ws.Range("H73").Value = 39052.5
$ws.Range("I73").Value = 38000
$ws.Range("K73").Value = 38000
$ws.Range("M73").Value = -36908
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 1244.3889
$ws.Range("I113").Value = 1075.4073
$ws.Range("K113").Value = 3226.2219
$ws.Range("M113").Value = -1056.2219
